# Group 2 MRI data layer2 - update subject IDs (SubjectID1..5 -> SubjectID6..10)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = "SubjectID6"
$ws.Range("A3").Value = "SubjectID7"
$ws.Range("A4").Value = "SubjectID8"
$ws.Range("A5").Value = "SubjectID9"
$ws.Range("A6").Value = "SubjectID10"

# Match the author's final selection state (B2:C6, active cell B2)
$ws.Range("B2:C6").Select()
